# Update latest output (run 7)
# Applies the "optimisation_result" refresh: new Schedule totals (E2/F2)
# and refreshed Detailed price history/forecast values for rows 7-49,
# including two rows whose Type flips from "forecast" to "historical".

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: updated cost / unit cost for the single schedule row ---
$schedule.Range("E2").Value = 1718.68822125
$schedule.Range("F2").Value = 28.41746397569445

# --- Detailed sheet: refreshed Price column (and Type where it changed) ---
$detailed.Range("B7").Value = 78
$detailed.Range("B8").Value = 78

$detailed.Range("B9").Value = 67.27928
$detailed.Range("C9").Value = "historical"

$detailed.Range("B10").Value = 65.84793999999999
$detailed.Range("C10").Value = "historical"

$detailed.Range("B11").Value = 61.54031
$detailed.Range("B12").Value = 62.8085

$detailed.Range("B14").Value = 80.02
$detailed.Range("B15").Value = 78.33004
$detailed.Range("B16").Value = 52.33926
$detailed.Range("B17").Value = 51.39572
$detailed.Range("B18").Value = 50.07721

$detailed.Range("B20").Value = 42.1835
$detailed.Range("B21").Value = 19.39377
$detailed.Range("B22").Value = 25.87642
$detailed.Range("B23").Value = 36.06

$detailed.Range("B25").Value = 22.07

$detailed.Range("B28").Value = 23.65517

$detailed.Range("B30").Value = 52.11771
$detailed.Range("B31").Value = 59.75743
$detailed.Range("B32").Value = 59.25835

$detailed.Range("B34").Value = 47.57833
$detailed.Range("B35").Value = 61.7683
$detailed.Range("B36").Value = 62.02801
$detailed.Range("B37").Value = 28.59972
$detailed.Range("B38").Value = 57.09237
$detailed.Range("B39").Value = 73.69302
$detailed.Range("B40").Value = 135.68796
$detailed.Range("B41").Value = 165.37731
$detailed.Range("B42").Value = 167.77584

$detailed.Range("B48").Value = 63.56007
$detailed.Range("B49").Value = 59.821
